$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 126
$ws.Range("I33").Value = 70
$ws.Range("J33").Value = 200.66667
$ws.Range("K33").Value = 70
$ws.Range("L33").Value = 200.66667
$ws.Range("M33").Value = 159
$ws.Range("N33").Value = -658.6666700000001

$ws.Range("H74").Value = 3325
$ws.Range("J74").Value = 3325
$ws.Range("L74").Value = 3325
$ws.Range("N74").Value = -5197

$ws.Range("H77").Value = 3325
$ws.Range("J77").Value = 3325
$ws.Range("L77").Value = 16625
$ws.Range("N77").Value = -25985

$ws.Range("H94").Value = 2483.889
$ws.Range("I94").Value = 2483.889
$ws.Range("K94").Value = 2483.889
$ws.Range("M94").Value = -2032.889

$ws.Range("H98").Value = 1610.1333
$ws.Range("I98").Value = 1582.2858
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 1582.2858
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -84.28580000000011
$ws.Range("N98").Value = -4996

$ws.Range("H100").Value = 1375.0834
$ws.Range("I100").Value = 1410.1
$ws.Range("J100").Value = 1200
$ws.Range("K100").Value = 1410.1
$ws.Range("L100").Value = 1200
$ws.Range("M100").Value = -869.0999999999999
$ws.Range("N100").Value = -2282

$ws.Range("H103").Value = 1069.4166
$ws.Range("J103").Value = 1076.6364
$ws.Range("L103").Value = 3229.9092
$ws.Range("N103").Value = -4401.9092

$ws.Range("H122").Value = 1610.1333
$ws.Range("I122").Value = 1582.2858
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4746.857400000001
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2296.857400000001
$ws.Range("N122").Value = -10900

$ws.Range("H125").Value = 1771.2222
$ws.Range("I125").Value = 583.3333
$ws.Range("J125").Value = 2008.8
$ws.Range("K125").Value = 5249.9997
$ws.Range("L125").Value = 18079.2
$ws.Range("M125").Value = -2789.9997
$ws.Range("N125").Value = -22999.2

$ws.Range("H126").Value = 57986.668
$ws.Range("J126").Value = 57986.668
$ws.Range("L126").Value = 57986.668
$ws.Range("N126").Value = -67866.66800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1464.3438
$ws.Range("I2").Value = 1332.238
$ws.Range("J2").Value = 1716.5454
$ws.Range("K2").Value = 1332.238
$ws.Range("L2").Value = 1716.5454
$ws.Range("M2").Value = -1219.238
$ws.Range("N2").Value = -1942.5454

$ws.Range("H45").Value = 1630.7142
$ws.Range("I45").Value = 1349.7894
$ws.Range("J45").Value = 2223.7778
$ws.Range("K45").Value = 1349.7894
$ws.Range("L45").Value = 2223.7778
$ws.Range("M45").Value = -972.7893999999999
$ws.Range("N45").Value = -2977.7778

$ws.Range("H61").Value = 4923.5713
$ws.Range("I61").Value = 1939.5
$ws.Range("J61").Value = 8902.333000000001
$ws.Range("K61").Value = 1939.5
$ws.Range("L61").Value = 8902.333000000001
$ws.Range("M61").Value = -1727.5
$ws.Range("N61").Value = -9326.333000000001

$ws.Range("H74").Value = 3357.878
$ws.Range("I74").Value = 4161.1665
$ws.Range("J74").Value = 1167.091
$ws.Range("K74").Value = 4161.1665
$ws.Range("L74").Value = 1167.091
$ws.Range("M74").Value = -3287.1665
$ws.Range("N74").Value = -2915.091

$ws.Range("H77").Value = 3357.878
$ws.Range("I77").Value = 4161.1665
$ws.Range("J77").Value = 1167.091
$ws.Range("K77").Value = 20805.8325
$ws.Range("L77").Value = 5835.455
$ws.Range("M77").Value = -16437.8325
$ws.Range("N77").Value = -14571.455

$ws.Range("H102").Value = 1966
$ws.Range("I102").Value = 1899.25
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 1899.25
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -277.25
$ws.Range("N102").Value = -5744

$ws.Range("H116").Value = 1464.3438
$ws.Range("I116").Value = 1332.238
$ws.Range("J116").Value = 1716.5454
$ws.Range("K116").Value = 1332.238
$ws.Range("L116").Value = 1716.5454
$ws.Range("M116").Value = 961.7619999999999
$ws.Range("N116").Value = -6304.5454

$ws.Range("H132").Value = 1462.34
$ws.Range("I132").Value = 1148.25
$ws.Range("J132").Value = 3765.6667
$ws.Range("K132").Value = 3444.75
$ws.Range("L132").Value = 11297.0001
$ws.Range("M132").Value = -914.75
$ws.Range("N132").Value = -16357.0001

$ws.Range("H136").Value = 4923.5713
$ws.Range("I136").Value = 1939.5
$ws.Range("J136").Value = 8902.333000000001
$ws.Range("K136").Value = 5818.5
$ws.Range("L136").Value = 26706.999
$ws.Range("M136").Value = -3268.5
$ws.Range("N136").Value = -31806.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1464.3438
$ws.Range("I3").Value = 1332.238
$ws.Range("J3").Value = 1716.5454
$ws.Range("K3").Value = 1332.238
$ws.Range("L3").Value = 1716.5454
$ws.Range("M3").Value = -1218.238
$ws.Range("N3").Value = -1944.5454

$ws.Range("H92").Value = 51238
$ws.Range("J92").Value = 51238
$ws.Range("L92").Value = 51238
$ws.Range("N92").Value = -56230

$ws.Range("H99").Value = 2577.875
$ws.Range("I99").Value = 1327.5
$ws.Range("J99").Value = 3828.25
$ws.Range("K99").Value = 1327.5
$ws.Range("L99").Value = 3828.25
$ws.Range("M99").Value = 170.5
$ws.Range("N99").Value = -6824.25

$ws.Range("H105").Value = 2626.41
$ws.Range("I105").Value = 1538.6875
$ws.Range("J105").Value = 2833.5952
$ws.Range("K105").Value = 1538.6875
$ws.Range("L105").Value = 2833.5952
$ws.Range("M105").Value = 208.3125
$ws.Range("N105").Value = -6327.5952

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 93.07692
$ws.Range("I7").Value = 88.5
$ws.Range("J7").Value = 100.4
$ws.Range("K7").Value = 88.5
$ws.Range("L7").Value = 100.4
$ws.Range("M7").Value = 24.5
$ws.Range("N7").Value = -326.4

$ws.Range("H16").Value = 1173.2222
$ws.Range("I16").Value = 1173.2222
$ws.Range("K16").Value = 1173.2222
$ws.Range("M16").Value = -886.2221999999999

$ws.Range("H31").Value = 1667.5625
$ws.Range("I31").Value = 1256.122
$ws.Range("J31").Value = 4077.4285
$ws.Range("K31").Value = 1256.122
$ws.Range("L31").Value = 4077.4285
$ws.Range("M31").Value = -961.1220000000001
$ws.Range("N31").Value = -4667.4285

$ws.Range("H34").Value = 1667.5625
$ws.Range("I34").Value = 1256.122
$ws.Range("J34").Value = 4077.4285
$ws.Range("K34").Value = 1256.122
$ws.Range("L34").Value = 4077.4285
$ws.Range("M34").Value = -1054.122
$ws.Range("N34").Value = -4481.4285

$ws.Range("H99").Value = 1578.5172
$ws.Range("I99").Value = 1323.1904
$ws.Range("K99").Value = 1323.1904
$ws.Range("M99").Value = 174.8096

$ws.Range("H107").Value = 610.29785
$ws.Range("I107").Value = 453.2069
$ws.Range("J107").Value = 863.3889
$ws.Range("K107").Value = 453.2069
$ws.Range("L107").Value = 863.3889
$ws.Range("M107").Value = 1466.7931
$ws.Range("N107").Value = -4703.3889

$ws.Range("H113").Value = 1173.2222
$ws.Range("I113").Value = 1173.2222
$ws.Range("K113").Value = 1173.2222
$ws.Range("M113").Value = 996.7778000000001

$ws.Range("H122").Value = 1395.4667
$ws.Range("I122").Value = 1142.2
$ws.Range("K122").Value = 3426.6
$ws.Range("M122").Value = -976.6000000000004

$ws.Range("H126").Value = 1578.5172
$ws.Range("I126").Value = 1323.1904
$ws.Range("K126").Value = 3969.5712
$ws.Range("M126").Value = -1499.5712

$ws.Range("H132").Value = 1995.2046
$ws.Range("I132").Value = 1269.1
$ws.Range("J132").Value = 3551.1428
$ws.Range("K132").Value = 3807.3
$ws.Range("L132").Value = 10653.4284
$ws.Range("M132").Value = -1277.3
$ws.Range("N132").Value = -15713.4284

$ws.Range("H134").Value = 1774.7073
$ws.Range("I134").Value = 1189.4445
$ws.Range("J134").Value = 2903.4285
$ws.Range("K134").Value = 3568.3335
$ws.Range("L134").Value = 8710.2855
$ws.Range("M134").Value = -1033.3335
$ws.Range("N134").Value = -13780.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1670.6389
$ws.Range("I102").Value = 1689.0741
$ws.Range("K102").Value = 1689.0741
$ws.Range("M102").Value = -67.07410000000004

$ws.Range("H113").Value = 7049.278
$ws.Range("I113").Value = 1002
$ws.Range("K113").Value = 1002
$ws.Range("M113").Value = 1168

$ws.Range("H123").Value = 19329.666
$ws.Range("J123").Value = 19329.666
$ws.Range("L123").Value = 19329.666
$ws.Range("N123").Value = -24229.666

$ws.Range("H132").Value = 2405.5334
$ws.Range("I132").Value = 2011.0435
$ws.Range("J132").Value = 3701.7144
$ws.Range("K132").Value = 6033.1305
$ws.Range("L132").Value = 11105.1432
$ws.Range("M132").Value = -3503.1305
$ws.Range("N132").Value = -16165.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1533.2646
$ws.Range("I7").Value = 1449.3448
$ws.Range("K7").Value = 1449.3448
$ws.Range("M7").Value = -1337.3448

$ws.Range("H61").Value = 8282.143
$ws.Range("I61").Value = 8282.143
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 8282.143
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -8080.143
$ws.Range("N61").ClearContents()

$ws.Range("H113").Value = 8282.143
$ws.Range("I113").Value = 8282.143
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 8282.143
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -6112.143
$ws.Range("N113").ClearContents()

$ws.Range("H126").Value = 1533.2646
$ws.Range("I126").Value = 1449.3448
$ws.Range("K126").Value = 4348.0344
$ws.Range("M126").Value = -1878.0344

$ws.Range("H132").Value = 7815.659
$ws.Range("I132").Value = 9550.3125
$ws.Range("J132").Value = 3189.9167
$ws.Range("K132").Value = 28650.9375
$ws.Range("L132").Value = 9569.750100000001
$ws.Range("M132").Value = -26120.9375
$ws.Range("N132").Value = -14629.7501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 4000
$ws.Range("I22").Value = 4000
$ws.Range("K22").Value = 4000
$ws.Range("M22").Value = -3707

$ws.Range("H122").Value = 2652.5356
$ws.Range("I122").Value = 2529.65
$ws.Range("K122").Value = 7588.950000000001
$ws.Range("M122").Value = -5138.950000000001
